$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests_AOS")

$ws.Range("J16").Value = "gal123473"
$ws.Range("J17").Value = "gal12l3@gmail.com"
$ws.Range("J18").Value = "Gal23416"
$ws.Range("J19").Value = "el12345"
$ws.Range("J20").Value = "El23456"

# K23 and K24 look numeric ("02", "2025") but must stay as text, matching
# the original inline-string cell type. Force text via NumberFormat, set
# the value, then clear the format again so no new style is left applied
# on the cell (keeps formatting identical to before the edit).
$k23 = $ws.Range("K23")
$k23.NumberFormat = "@"
$k23.Value = "02"
$k23.ClearFormats()

$k24 = $ws.Range("K24")
$k24.NumberFormat = "@"
$k24.Value = "2025"
$k24.ClearFormats()

$ws.Range("K25").Value = "gal-elad"
